$wb = $excel.ActiveWorkbook

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 56049.5
$ws.Range("I6").Value = 215.5
$ws.Range("K6").Value = 646.5
$ws.Range("M6").Value = -534.5

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 397.44
$ws.Range("I19").Value = 217.22223
$ws.Range("J19").Value = 498.8125
$ws.Range("K19").Value = 217.22223
$ws.Range("L19").Value = 498.8125
$ws.Range("M19").Value = -42.22223
$ws.Range("N19").Value = -848.8125

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2534
$ws.Range("I40").Value = 3184.5
$ws.Range("J40").Value = 1666.6666
$ws.Range("K40").Value = 3184.5
$ws.Range("L40").Value = 1666.6666
$ws.Range("M40").Value = -3009.5
$ws.Range("N40").Value = -2016.6666

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4318.65
$ws.Range("I64").Value = 3914.2856
$ws.Range("J64").Value = 4536.385
$ws.Range("K64").Value = 3914.2856
$ws.Range("L64").Value = 4536.385
$ws.Range("M64").Value = -3666.2856
$ws.Range("N64").Value = -5032.385

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4318.65
$ws.Range("I67").Value = 3914.2856
$ws.Range("J67").Value = 4536.385
$ws.Range("K67").Value = 3914.2856
$ws.Range("L67").Value = 4536.385
$ws.Range("M67").Value = -3056.2856
$ws.Range("N67").Value = -6252.385

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2001386.1
$ws.Range("I137").Value = 3572370.2
$ws.Range("J137").Value = 1951.909
$ws.Range("K137").Value = 10717110.6
$ws.Range("L137").Value = 5855.727000000001
$ws.Range("M137").Value = -10714560.6
$ws.Range("N137").Value = -10955.727

# ARM sheet updates
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 7479.5713
$ws.Range("I88").Value = 3195.8
$ws.Range("J88").Value = 9859.444
$ws.Range("K88").Value = 3195.8
$ws.Range("L88").Value = 9859.444
$ws.Range("M88").Value = -2789.8
$ws.Range("N88").Value = -10671.444

# ARM sheet updates
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 7479.5713
$ws.Range("I91").Value = 3195.8
$ws.Range("J91").Value = 9859.444
$ws.Range("K91").Value = 3195.8
$ws.Range("L91").Value = 9859.444
$ws.Range("M91").Value = -1791.8
$ws.Range("N91").Value = -12667.444

# ARM sheet updates
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 70155
$ws.Range("J106").Value = 70155
$ws.Range("L106").Value = 70155
$ws.Range("N106").Value = -72679

# ARM sheet updates
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 89719.46000000001
$ws.Range("I132").Value = 61074.59
$ws.Range("K132").Value = 183223.77
$ws.Range("M132").Value = -180693.77

# BSM sheet updates
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17399.762
$ws.Range("I86").Value = 17842.215
$ws.Range("J86").Value = 16514.857
$ws.Range("K86").Value = 17842.215
$ws.Range("L86").Value = 16514.857
$ws.Range("M86").Value = -16719.215
$ws.Range("N86").Value = -18760.857

# BSM sheet updates
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 17399.762
$ws.Range("I89").Value = 17842.215
$ws.Range("J89").Value = 16514.857
$ws.Range("K89").Value = 89211.075
$ws.Range("L89").Value = 82574.285
$ws.Range("M89").Value = -83595.075
$ws.Range("N89").Value = -93806.285

# BSM sheet updates
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2014.6666
$ws.Range("I134").Value = 2349.8462
$ws.Range("K134").Value = 7049.5386
$ws.Range("M134").Value = -4514.5386

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 455735.7
$ws.Range("I68").Value = 899
$ws.Range("J68").Value = 1430385.8
$ws.Range("K68").Value = 2697
$ws.Range("L68").Value = 4291157.4
$ws.Range("M68").Value = -1886
$ws.Range("N68").Value = -4292779.4

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 455735.7
$ws.Range("I71").Value = 899
$ws.Range("J71").Value = 1430385.8
$ws.Range("K71").Value = 8091
$ws.Range("L71").Value = 12873472.2
$ws.Range("M71").Value = -4035
$ws.Range("N71").Value = -12881584.2

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2056.2974
$ws.Range("I80").Value = 1452.3077
$ws.Range("J80").Value = 2383.4583
$ws.Range("K80").Value = 4356.9231
$ws.Range("L80").Value = 7150.374899999999
$ws.Range("M80").Value = -3420.9231
$ws.Range("N80").Value = -9022.374899999999

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2056.2974
$ws.Range("I83").Value = 1452.3077
$ws.Range("J83").Value = 2383.4583
$ws.Range("K83").Value = 13070.7693
$ws.Range("L83").Value = 21451.1247
$ws.Range("M83").Value = -8390.7693
$ws.Range("N83").Value = -30811.1247

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 552
$ws.Range("I113").Value = 490.4054
$ws.Range("J113").Value = 686.05884
$ws.Range("K113").Value = 1471.2162
$ws.Range("L113").Value = 2058.17652
$ws.Range("M113").Value = 698.7838000000002
$ws.Range("N113").Value = -6398.17652

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1239.2388
$ws.Range("J131").Value = 1347.4822
$ws.Range("L131").Value = 4042.4466
$ws.Range("N131").Value = -14122.4466

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5000000
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5000000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4999888
$ws.Range("N7").ClearContents()

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 5000000
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -4999861
$ws.Range("N8").ClearContents()

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30590.41
$ws.Range("I70").Value = 41722.89
$ws.Range("J70").Value = 5542.3335
$ws.Range("K70").Value = 41722.89
$ws.Range("L70").Value = 5542.3335
$ws.Range("M70").Value = -41452.89
$ws.Range("N70").Value = -6082.3335

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 30590.41
$ws.Range("I73").Value = 41722.89
$ws.Range("J73").Value = 5542.3335
$ws.Range("K73").Value = 41722.89
$ws.Range("L73").Value = 5542.3335
$ws.Range("M73").Value = -40786.89
$ws.Range("N73").Value = -7414.3335

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 67317.44500000001
$ws.Range("J101").Value = 67317.44500000001
$ws.Range("L101").Value = 67317.44500000001
$ws.Range("N101").Value = -73807.44500000001

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 96835.19
$ws.Range("I132").Value = 63596.938
$ws.Range("J132").Value = 203197.6
$ws.Range("K132").Value = 190790.814
$ws.Range("L132").Value = 609592.8
$ws.Range("M132").Value = -188260.814
$ws.Range("N132").Value = -614652.8

# WVR sheet updates
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# WVR sheet updates
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4867.3335
$ws.Range("I14").Value = 204
$ws.Range("J14").Value = 5800
$ws.Range("K14").Value = 204
$ws.Range("L14").Value = 5800
$ws.Range("M14").Value = -36
$ws.Range("N14").Value = -6136

# WVR sheet updates
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 37284.6
$ws.Range("I136").Value = 20546.137
$ws.Range("J136").Value = 250700
$ws.Range("K136").Value = 61638.41099999999
$ws.Range("L136").Value = 752100
$ws.Range("M136").Value = -59088.41099999999
$ws.Range("N136").Value = -757200
